# Fix named-range requirement text: hyphens aren't allowed in Excel
# defined-name identifiers, so swap them for underscores / drop them.
#
#   Sheet1!A3   "Name this range "Nums1-10":"        -> "Nums1_10"
#   Sheet1!A13  "Sum of Nums1-10:"                    -> "Nums1_10"
#   Sheet1!A14  "Average of Nums1-10:"                -> "Nums1_10"
#   Sheet2!A27  "1. Create ... "Total-Sales"."        -> "TotalSales"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A27").Value = '1. Create a separate Excel file called "TotalSales".'

$ws1.Range("A3").Value  = 'Name this range "Nums1_10":'
$ws1.Range("A13").Value = "Sum of Nums1_10:"
$ws1.Range("A14").Value = "Average of Nums1_10:"

# Re-point the on-screen selections to where the author left off while
# making the fix.
$ws2.Select()
$ws2.Range("A27").Select()

$ws1.Select()
$ws1.Range("A15:A19").Select()

# Restore the workbook window size recorded in the saved view state.
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 12240
